$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 4208.293336037124
$ws.Range("D2").Value = 1051.36775377502
